$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Oct 25 12:06:23 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 12:06:38 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 12:06:53 EDT 2024"
